$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '65.385.73'
Set-TextValue $ws.Range('E2') '  -3.34%  '

Set-TextValue $ws.Range('D3') '3.483.68'
Set-TextValue $ws.Range('E3') '  -0.18%  '

Set-TextValue $ws.Range('E4') '  +0.13%  '

Set-TextValue $ws.Range('D5') '552.60'
Set-TextValue $ws.Range('E5') '  -0.24%  '

Set-TextValue $ws.Range('D6') '178.02'
Set-TextValue $ws.Range('E6') '  -5.21%  '

Set-TextValue $ws.Range('E7') '  +4.26%  '

Set-TextValue $ws.Range('E8') '  +0.04%  '

Set-TextValue $ws.Range('E9') '  -1.17%  '

Set-TextValue $ws.Range('E10') '  +1.87%  '

Set-TextValue $ws.Range('E11') '  -5.68%  '

Set-TextValue $ws.Range('E12') '  -1.90%  '

Set-TextValue $ws.Range('D13') '9.21'
Set-TextValue $ws.Range('E13') '  -2.73%  '

Set-TextValue $ws.Range('D14') '4.044.30'
Set-TextValue $ws.Range('E14') '  -0.28%  '

Set-TextValue $ws.Range('D15') '18.51'
Set-TextValue $ws.Range('E15') '  +1.18%  '

Set-TextValue $ws.Range('D16') '3.481.83'
Set-TextValue $ws.Range('E16') '  -0.30%  '

Set-TextValue $ws.Range('E17') '  +0.48%  '

Set-TextValue $ws.Range('D18') '12.04'
Set-TextValue $ws.Range('E18') '  +1.82%  '

Set-TextValue $ws.Range('D19') '65.440.42'
Set-TextValue $ws.Range('E19') '  -3.58%  '

Set-TextValue $ws.Range('D20') '0.987'
Set-TextValue $ws.Range('E20') '  -2.30%  '

Set-TextValue $ws.Range('D21') '415.89'
Set-TextValue $ws.Range('E21') '  +2.79%  '

Set-TextValue $ws.Range('E22') '  +1.72%  '

Set-TextValue $ws.Range('D23') '86.13'
Set-TextValue $ws.Range('E23') '  +1.60%  '

Set-TextValue $ws.Range('D24') '4.26'
Set-TextValue $ws.Range('E24') '  +1.20%  '

Set-TextValue $ws.Range('D25') '12.74'
Set-TextValue $ws.Range('E25') '  +7.47%  '

Set-TextValue $ws.Range('E26') '  -11.35%  '

Set-TextValue $ws.Range('D27') '2.84'
Set-TextValue $ws.Range('E27') '  -3.47%  '

Set-TextValue $ws.Range('D28') '6.03'
Set-TextValue $ws.Range('E28') '  -3.50%  '

Set-TextValue $ws.Range('D29') '9.00'
Set-TextValue $ws.Range('E29') '  +4.18%  '

Set-TextValue $ws.Range('D30') '30.16'

Set-TextValue $ws.Range('D31') '6.48'
Set-TextValue $ws.Range('E31') '  -6.39%  '

Set-TextValue $ws.Range('D32') '609.30'
Set-TextValue $ws.Range('E32') '  -10.96%  '

Set-TextValue $ws.Range('D33') '11.69'

Set-TextValue $ws.Range('E34') '  -1.37%  '

Set-TextValue $ws.Range('D35') '59.46'
Set-TextValue $ws.Range('E35') '  -0.36%  '

Set-TextValue $ws.Range('E36') '  +9.42%  '

Set-TextValue $ws.Range('D37') '1.00'
Set-TextValue $ws.Range('E37') '  +0.29%  '

Set-TextValue $ws.Range('D38') '37.25'
Set-TextValue $ws.Range('E38') '  -4.26%  '

Set-TextValue $ws.Range('B39') 'Maker'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D39') '3.368.29'
Set-TextValue $ws.Range('E39') '  +10.05%  '

Set-TextValue $ws.Range('B40') 'PEPE'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D40') '0.0₃0782'
Set-TextValue $ws.Range('E40') '  -6.32%  '

Set-TextValue $ws.Range('E41') '  -5.92%  '

Set-TextValue $ws.Range('D42') '0.998'
Set-TextValue $ws.Range('E42') '  -0.47%  '

Set-TextValue $ws.Range('E43') '  -6.99%  '

Set-TextValue $ws.Range('D44') '2.83'
Set-TextValue $ws.Range('E44') '  -4.42%  '

Set-TextValue $ws.Range('E45') '  -10.30%  '

Set-TextValue $ws.Range('E46') '  -2.33%  '

Set-TextValue $ws.Range('E47') '  -3.39%  '

Set-TextValue $ws.Range('E48') '  -1.25%  '

Set-TextValue $ws.Range('E49') '  +1.28%  '

Set-TextValue $ws.Range('E50') '  -6.06%  '

Set-TextValue $ws.Range('D51') '137.89'
Set-TextValue $ws.Range('E51') '  -1.33%  '
